$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. "1.004", "1.000").
# Mark each changed Price cell as Text *before* writing its new value so
# Excel keeps the exact literal string instead of auto-converting it to a
# number (which would also silently drop meaningful trailing zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.353.40'
$ws.Range("E2").Value = '  +6.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.816.40'
$ws.Range("E3").Value = '  +6.33%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  +0.48%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '344.82'
$ws.Range("E5").Value = '  +3.92%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  +0.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3848'
$ws.Range("E7").Value = '  +4.32%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.37'
$ws.Range("E8").Value = '  +3.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3531'
$ws.Range("E9").Value = '  +6.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.244'
$ws.Range("E10").Value = '  +6.41%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07798'
$ws.Range("E11").Value = '  +6.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.75'
$ws.Range("E12").Value = '  +13.86%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.003'
$ws.Range("E13").Value = '  +0.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.673'
$ws.Range("E14").Value = '  +7.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.261'
$ws.Range("E15").Value = '  +5.51%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.815.64'
$ws.Range("E16").Value = '  +6.68%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001133'
$ws.Range("E17").Value = '  +5.77%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06779'
$ws.Range("E18").Value = '  +2.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.34'
$ws.Range("E19").Value = '  +7.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.40%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.98'
$ws.Range("E21").Value = '  +10.90%  '

$ws.Range("E22").Value = '  +8.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.18'
$ws.Range("E23").Value = '  +1.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.356.82'
$ws.Range("E24").Value = '  +6.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.468'
$ws.Range("E25").Value = '  +0.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.748'
$ws.Range("E26").Value = '  +9.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '22.09'
$ws.Range("E27").Value = '  +15.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '1.509'
$ws.Range("E28").Value = '  +16.63%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.33'
$ws.Range("E29").Value = '  +3.20%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.019.05'
$ws.Range("E30").Value = '  +6.78%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '137.53'
$ws.Range("E31").Value = '  +7.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.433'
$ws.Range("E32").Value = '  +7.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.123'
$ws.Range("E33").Value = '  +0.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.90'
$ws.Range("E34").Value = '  +8.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08849'
$ws.Range("E35").Value = '  +4.06%  '

$ws.Range("E36").Value = '  +3.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '5.672'
$ws.Range("E37").Value = '  +6.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").Value = '0.7178'
$ws.Range("E38").Value = '  +17.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06595'
$ws.Range("E39").Value = '  +5.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2285'
$ws.Range("E40").Value = '  +7.66%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02431'
$ws.Range("E41").Value = '  +7.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.096'
$ws.Range("E42").Value = '  +6.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.262'
$ws.Range("E43").Value = '  -0.86%  '

$ws.Range("E44").Value = '  +4.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6721'
$ws.Range("E45").Value = '  +14.96%  '

$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.042'
$ws.Range("E47").Value = '  +4.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.199'
$ws.Range("E48").Value = '  +9.52%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '133.55'
$ws.Range("E49").Value = '  +5.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07377'
$ws.Range("E50").Value = '  +2.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.24'
$ws.Range("E51").Value = '  +5.76%  '
